$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.838.21'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.483.98'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -6.05%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '557.98'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.37'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -5.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.603'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.483.55'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -5.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.109'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -8.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.51'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -5.28%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.359'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -6.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.65'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -7.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.927.55'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -6.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000169'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -8.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.712.15'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.479.05'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -6.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.24'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -7.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.25'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -7.52%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -6.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '323.12'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -6.43%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.91'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '64.59'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -5.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000101'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -9.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '567.02'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.608.22'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -5.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.52'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -7.84%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.42'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -10.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.80'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.151'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -6.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.95'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.62'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -6.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.03'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -8.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.98'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -9.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.998'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.64'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -5.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '144.98'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.48%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -6.97%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.49'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '40.71'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '149.34'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -8.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '22.28'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -8.62%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.66'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -6.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0546'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -7.60%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.602'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -5.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0946'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -5.47%  '
